$d = $word.ActiveDocument

# --- Simple text replacements -------------------------------------------------
$d.Content.Find.Execute("564564654", $false, $false, $false, $false, $false, $true, 1, $false, "05465465465", 2) | Out-Null
$d.Content.Find.Execute("dgsdgdsgdsg", $false, $false, $false, $false, $false, $true, 1, $false, "hassanalihazaraa@gmail.com", 2) | Out-Null
$d.Content.Find.Execute("dsgsdgsdg", $false, $false, $false, $false, $false, $true, 1, $false, "Hello world", 2) | Out-Null
$d.Content.Find.Execute("sdgsdgsdg", $false, $false, $false, $false, $false, $true, 1, $false, "Becode ", 2) | Out-Null
$d.Content.Find.Execute("2000-2020", $false, $false, $false, $false, $false, $true, 1, $false, "2020-2020", 2) | Out-Null

# --- Append a new "Risk concile" work-experience paragraph -------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$newPara = $lastPara.Range.InsertParagraphAfter()
$d.Paragraphs($d.Paragraphs.Count - 1).Range.Text = ""

$bRange = $d.Paragraphs($d.Paragraphs.Count - 1).Range
$bRange.Text = "Risk concile "
$bRange.Font.Bold = $true
$bRange.Font.Italic = $false

$endOfPara = $d.Paragraphs($d.Paragraphs.Count - 1).Range
$endOfPara.Collapse(0)
$endOfPara.InsertAfter("2020-2025")
$endOfPara.Font.Bold = $false
$endOfPara.Font.Italic = $true
$endOfPara.Collapse(0)
$endOfPara.InsertAfter("`v")

# --- Skills heading + bullet list ---------------------------------------------
$skillsPara = $d.Paragraphs($d.Paragraphs.Count).Range.InsertParagraphAfter()
$skillsRange = $d.Paragraphs($d.Paragraphs.Count).Range
$skillsRange.Text = "Skills"
$skillsRange.Style = "Heading1"

$skills = @("HTML", "CSS", "Javascript", "PHP", "Symfony", "Mysql", "Python", "Django", "Java", "Spring boot")
foreach ($skill in $skills) {
    $p = $d.Paragraphs($d.Paragraphs.Count).Range.InsertParagraphAfter()
    $r = $d.Paragraphs($d.Paragraphs.Count).Range
    $r.Text = $skill
    $r.Style = "List Bullet"
}

# --- Footer --------------------------------------------------------------------
$footer = $d.Sections(1).Footers(1)
$footer.Range.Text = "CV generated using Hassan`"s app"
$footer.Range.Style = "Footer"
